$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the runs of the "...burning in the sunlight for example)."
#    sentence (removing the gramStart/gramEnd proofErr markers around
#    "sunlight") without also merging it into the preceding, unrelated
#    run ("find a new functional body...").
# ------------------------------------------------------------------
$soulPara = $d.Paragraphs.Item(5)
$soulFull = $soulPara.Range.Text
$soulStart = $soulPara.Range.Start

# Temporarily give the preceding run distinct formatting so the
# upcoming Find/Replace (which merges same-formatted adjacent runs)
# does not sweep it into the replacement.
$bodyRunStart = $soulStart + $soulFull.IndexOf("find a new functional body")
$bodyRunEnd = $soulStart + $soulFull.IndexOf("That comes with all the advantages")
$bodyRun = $d.Range($bodyRunStart, $bodyRunEnd)
$bodyRun.Bold = 1

$sentenceRange = $soulPara.Range.Duplicate
$sentenceRange.Find.Execute(
    "the sunlight for example).", $true, $false, $false, $false, $false,
    $true, 1, $false, "the sunlight for example).", 2)

# Restore the preceding run's formatting (this is a self-contained
# formatting-only edit, so it will not re-merge the runs).
$bodyRun2Start = $bodyRunStart
$bodyRun2End = $bodyRunEnd
$bodyRun2 = $d.Range($bodyRun2Start, $bodyRun2End)
$bodyRun2.Bold = 0

# ------------------------------------------------------------------
# 2) Insert a new list paragraph (same list as the surrounding ones)
#    describing the new "Rebirth" potion effect.
# ------------------------------------------------------------------
$soulPara = $d.Paragraphs.Item(5)
$soulPara.Range.InsertParagraphAfter()
$newParaIndex = 6

$segments = @(
    "Potion effect, Rebirth: Gradually transforms a lich player into a",
    " living player",
    " (heart by heart)",
    ", at the rate of ",
    "1 heart per second.",
    " "
)
$segmentsText = [string]::Join("", $segments)

$newPara = $d.Paragraphs.Item($newParaIndex)
$insertPos = $newPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter($segmentsText)

# Split the freshly inserted text back into separate runs (matching
# how the original edit was authored) by toggling formatting on/off
# across each segment boundary.
$cursor = $insertPos
foreach ($segment in $segments) {
    $segStart = $cursor
    $segEnd = $cursor + $segment.Length
    $segRange = $d.Range($segStart, $segEnd)
    $segRange.Bold = 1
    $segRange.Bold = 0
    $cursor = $segEnd
}

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from its old location (after "you
#    have a strong soul") to right after "1 heart per second." (i.e.
#    just before the trailing space run) - mirroring how Word tracks
#    the last edit point. Adding a bookmark with an existing name
#    moves it, so the old one disappears automatically.
# ------------------------------------------------------------------
$bookmarkPos = $insertPos + ($segmentsText.Length - 1)
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
